$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: the phone number in A15 was entered as text; convert it to a
# genuine numeric value (matches the other phone columns in the sheet).
$ws.Range("A15").Value = 71277620

# Row 16: new redemption record for phone 71277620, 760 points redeemed.
# Phone numbers in this sheet are stored as text, so force text storage
# for the numeric-looking value (leading apostrophe = "treat as text"),
# then strip the format Excel applies for that so no stray style sticks.
$ws.Range("A16").Value = "'71277620"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = 760
$ws.Range("C16").Value = "2025-08-18T17:04:26"
